# Update crypto price/volume figures (and swap Polkadot/WrappedEther row order)
# per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.595.86"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "1.857.94"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'244.93"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.6932"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07692"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'0.3062"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'23.70"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "'0.07767"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.862.98"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.142"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "'90.94"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'0.6917"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "'6.563"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "29.475.41"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "'0.000008297"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "2.101.84"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "'239.78"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'12.76"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'7.605"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'0.1496"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").Value = "'8.915"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "'159.34"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "'18.27"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "'1.538"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "'4.249"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'4.172"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "'1.205"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").Value = "'0.7695"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "'1.891"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "'1.152"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "1.332.64"
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("D39").Value = "'0.01872"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").Value = "'2.728"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "'0.9689"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("D42").Value = "'106.49"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "'5.818"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").Value = "'9.780"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").Value = "1.999.72"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "'63.19"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'6.961"
$ws.Range("E51").Value = "  +0.80%  "
